$wb = $excel.ActiveWorkbook

# The workbook currently has a "事業投資" sheet holding investment data.
# We need to:
#   1) keep that data+formatting intact on a sheet literally named "事業投資"
#      (placed AFTER the new sheet, as sheetId 8)
#   2) introduce a new "債權" sheet (taking over the original sheet's slot /
#      sheetId 7) holding brand-new claims/debt data, re-using the same
#      header formatting as the original sheet.

$invest = $wb.Worksheets.Item("事業投資")

# Duplicate the investment sheet right after itself; the copy keeps all
# values + formatting (and becomes the tab that keeps the "事業投資" name).
$invest.Copy($null, $invest)
$investCopy = $wb.Worksheets.Item("事業投資 (2)")

# Turn the original sheet into "債權" first (freeing up the "事業投資" name)
# and swap in the new claims data, then rename the copy back to "事業投資".
$claims = $invest
$claims.Name = "債權"
$investCopy.Name = "事業投資"

# Wipe the old investment rows/headers out of the renamed sheet.
$claims.Range("A1:N5").ClearContents()

# Header row (row 1) - species/owner/debtor/total/... mirrors the layout
# used by every other property-type sheet in this workbook.
$claims.Range("B1").Value = "species"
$claims.Range("C1").Value = "owner"
$claims.Range("D1").Value = "debtor"
$claims.Range("E1").Value = "total"
$claims.Range("F1").Value = "register_date"
$claims.Range("G1").Value = "register_reason"
$claims.Range("H1").Value = "property_category"
$claims.Range("I1").Value = "category"
$claims.Range("J1").Value = "date"
$claims.Range("K1").Value = "legislator_name"
$claims.Range("L1").Value = "legislator_id"
$claims.Range("M1").Value = "source_file"
$claims.Range("N1").Value = "index"

# Re-apply the bold/centered/bordered header formatting used throughout the
# workbook by copying the format from an already-styled header cell.
$wb.Worksheets.Item("保險").Range("B1").Copy() | Out-Null
$claims.Range("A1:N1").PasteSpecial(-4122)

# Data rows 2-5.
$rows = @(
    @(297, "未兌現支票", "張慶忠", "漢寶開發建設(股)公司新北市中和區安樂路", 449020000, "98年10月23日", "借款", "claim", "normal", "2012-04-19", "張慶忠", 1347, "tmp93201", 297),
    @(298, "未兌現支票", "張慶忠", "資信建設有限公司新北市中和區安樂路",     268020000, "98年12月06日", "借款", "claim", "normal", "2012-04-19", "張慶忠", 1347, "tmp93201", 298),
    @(299, "未兌現支票", "張慶忠", "漢龍營造股份有限公司新北市中和區安樂路", 146800000, "98年11月23日", "借款", "claim", "normal", "2012-04-19", "張慶忠", 1347, "tmp93201", 299),
    @(300, "未兌現支票", "張慶忠", "芳林建設有限公司新北市中和區安樂路",     1000000,   "100年08月23日", "借款", "claim", "normal", "2012-04-19", "張慶忠", 1347, "tmp93201", 300)
)

$r = 2
foreach ($row in $rows) {
    $claims.Cells.Item($r, 1).Value = $row[0]
    $claims.Cells.Item($r, 2).Value = $row[1]
    $claims.Cells.Item($r, 3).Value = $row[2]
    $claims.Cells.Item($r, 4).Value = $row[3]
    $claims.Cells.Item($r, 5).Value = $row[4]
    $claims.Cells.Item($r, 6).Value = $row[5]
    $claims.Cells.Item($r, 7).Value = $row[6]
    $claims.Cells.Item($r, 8).Value = $row[7]
    $claims.Cells.Item($r, 9).Value = $row[8]
    $claims.Cells.Item($r, 10).Value = $row[9]
    $claims.Cells.Item($r, 11).Value = $row[10]
    $claims.Cells.Item($r, 12).Value = $row[11]
    $claims.Cells.Item($r, 13).Value = $row[12]
    $claims.Cells.Item($r, 14).Value = $row[13]
    $r = $r + 1
}

Write-Output "done"
